$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K6").Value = 8.200555000000001
$ws.Range("L6").Value = 8.200555000000001
$ws.Range("M6").Value = 9.480675000000002
$ws.Range("N6").Value = 9.480675000000002
$ws.Range("O6").Value = 10.408553
$ws.Range("P6").Value = 10.408553
$ws.Range("Q6").Value = 11.371678
$ws.Range("R6").Value = 12.516998
$ws.Range("S6").Value = 12.516998
